# Applies the weekly data refresh for Hortaliza / Macroferia Regional de Talca - Berenjena.
# A new observation (row 86) is inserted and every subsequent row shifts down by one,
# with the former last row (128) becoming the new last row (129).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 86
$ws.Range("D86").Value = 44726
$ws.Range("J86").Value = 200
$ws.Range("K86").Value = 8000
$ws.Range("L86").Value = 8000
$ws.Range("M86").Value = 8000
$ws.Range("N86").Value = '$/caja 50 unidades'
$ws.Range("O86").Value = 'Región de Arica y Parinacota'
$ws.Range("P86").Value = 160
$ws.Range("Q86").Value = 50

# Row 87
$ws.Range("D87").Value = 44575
$ws.Range("J87").Value = 150
$ws.Range("K87").Value = 8000
$ws.Range("L87").Value = 8000
$ws.Range("M87").Value = 8000
$ws.Range("N87").Value = '$/caja 60 unidades'
$ws.Range("O87").Value = 'Región del Maule'
$ws.Range("P87").Value = 133
$ws.Range("Q87").Value = 60

# Row 88
$ws.Range("D88").Value = 44257
$ws.Range("J88").Value = 150
$ws.Range("K88").Value = 8000
$ws.Range("L88").Value = 8000
$ws.Range("M88").Value = 8000
$ws.Range("N88").Value = '$/caja 60 unidades'
$ws.Range("O88").Value = 'Región del Maule'
$ws.Range("P88").Value = 133
$ws.Range("Q88").Value = 60

# Row 89
$ws.Range("D89").Value = 44596
$ws.Range("J89").Value = 150
$ws.Range("K89").Value = 7000
$ws.Range("L89").Value = 7000
$ws.Range("M89").Value = 7000
$ws.Range("N89").Value = '$/caja 50 unidades'
$ws.Range("O89").Value = 'Región del Maule'
$ws.Range("P89").Value = 140
$ws.Range("Q89").Value = 50

# Row 90
$ws.Range("D90").Value = 44467
$ws.Range("J90").Value = 300
$ws.Range("K90").Value = 7000
$ws.Range("L90").Value = 7000
$ws.Range("M90").Value = 7000
$ws.Range("N90").Value = '$/caja 60 unidades'
$ws.Range("O90").Value = 'Región de Arica y Parinacota'
$ws.Range("P90").Value = 117
$ws.Range("Q90").Value = 60

# Row 91
$ws.Range("D91").Value = 44448
$ws.Range("J91").Value = 200
$ws.Range("K91").Value = 7000
$ws.Range("L91").Value = 7000
$ws.Range("M91").Value = 7000
$ws.Range("N91").Value = '$/caja 50 unidades'
$ws.Range("O91").Value = 'Región de Arica y Parinacota'
$ws.Range("P91").Value = 140
$ws.Range("Q91").Value = 50

# Row 92
$ws.Range("D92").Value = 44435
$ws.Range("J92").Value = 300
$ws.Range("K92").Value = 6000
$ws.Range("L92").Value = 6000
$ws.Range("M92").Value = 6000
$ws.Range("N92").Value = '$/caja 50 unidades'
$ws.Range("O92").Value = 'Región de Arica y Parinacota'
$ws.Range("P92").Value = 120
$ws.Range("Q92").Value = 50

# Row 93
$ws.Range("D93").Value = 44435
$ws.Range("J93").Value = 1300
$ws.Range("K93").Value = 6000
$ws.Range("L93").Value = 7000
$ws.Range("M93").Value = 6231
$ws.Range("N93").Value = '$/caja 60 unidades'
$ws.Range("O93").Value = 'Región de Arica y Parinacota'
$ws.Range("P93").Value = 104
$ws.Range("Q93").Value = 60

# Row 94
$ws.Range("D94").Value = 44608
$ws.Range("J94").Value = 200
$ws.Range("K94").Value = 7000
$ws.Range("L94").Value = 7000
$ws.Range("M94").Value = 7000
$ws.Range("N94").Value = '$/caja 60 unidades'
$ws.Range("O94").Value = 'Región del Maule'
$ws.Range("P94").Value = 117
$ws.Range("Q94").Value = 60

# Row 95
$ws.Range("D95").Value = 44259
$ws.Range("J95").Value = 250
$ws.Range("K95").Value = 8000
$ws.Range("L95").Value = 8000
$ws.Range("M95").Value = 8000
$ws.Range("N95").Value = '$/caja 60 unidades'
$ws.Range("O95").Value = 'Región del Maule'
$ws.Range("P95").Value = 133
$ws.Range("Q95").Value = 60

# Row 96
$ws.Range("D96").Value = 44571
$ws.Range("J96").Value = 200
$ws.Range("K96").Value = 8000
$ws.Range("L96").Value = 8000
$ws.Range("M96").Value = 8000
$ws.Range("N96").Value = '$/caja 60 unidades'
$ws.Range("O96").Value = 'Provincia de Chacabuco'
$ws.Range("P96").Value = 133
$ws.Range("Q96").Value = 60

# Row 97
$ws.Range("D97").Value = 44418
$ws.Range("J97").Value = 200
$ws.Range("K97").Value = 8000
$ws.Range("L97").Value = 8000
$ws.Range("M97").Value = 8000
$ws.Range("N97").Value = '$/caja 60 unidades'
$ws.Range("O97").Value = 'Región de Arica y Parinacota'
$ws.Range("P97").Value = 133
$ws.Range("Q97").Value = 60

# Row 98
$ws.Range("D98").Value = 44284
$ws.Range("J98").Value = 200
$ws.Range("K98").Value = 8000
$ws.Range("L98").Value = 8000
$ws.Range("M98").Value = 8000
$ws.Range("N98").Value = '$/caja 60 unidades'
$ws.Range("O98").Value = 'Región del Maule'
$ws.Range("P98").Value = 133
$ws.Range("Q98").Value = 60

# Row 99
$ws.Range("D99").Value = 44663
$ws.Range("J99").Value = 150
$ws.Range("K99").Value = 10000
$ws.Range("L99").Value = 10000
$ws.Range("M99").Value = 10000
$ws.Range("N99").Value = '$/caja 50 unidades'
$ws.Range("O99").Value = 'Región del Maule'
$ws.Range("P99").Value = 200
$ws.Range("Q99").Value = 50

# Row 100
$ws.Range("D100").Value = 44424
$ws.Range("J100").Value = 300
$ws.Range("K100").Value = 8000
$ws.Range("L100").Value = 8000
$ws.Range("M100").Value = 8000
$ws.Range("N100").Value = '$/caja 60 unidades'
$ws.Range("O100").Value = 'Región de Arica y Parinacota'
$ws.Range("P100").Value = 133
$ws.Range("Q100").Value = 60

# Row 101
$ws.Range("D101").Value = 44722
$ws.Range("J101").Value = 300
$ws.Range("K101").Value = 6000
$ws.Range("L101").Value = 6000
$ws.Range("M101").Value = 6000
$ws.Range("N101").Value = '$/caja 50 unidades'
$ws.Range("O101").Value = 'Región de Arica y Parinacota'
$ws.Range("P101").Value = 120
$ws.Range("Q101").Value = 50

# Row 102
$ws.Range("D102").Value = 44664
$ws.Range("J102").Value = 200
$ws.Range("K102").Value = 9000
$ws.Range("L102").Value = 9000
$ws.Range("M102").Value = 9000
$ws.Range("N102").Value = '$/caja 50 unidades'
$ws.Range("O102").Value = 'Región del Maule'
$ws.Range("P102").Value = 180
$ws.Range("Q102").Value = 50

# Row 103
$ws.Range("D103").Value = 44274
$ws.Range("J103").Value = 150
$ws.Range("K103").Value = 7000
$ws.Range("L103").Value = 7000
$ws.Range("M103").Value = 7000
$ws.Range("N103").Value = '$/caja 60 unidades'
$ws.Range("O103").Value = 'Región del Maule'
$ws.Range("P103").Value = 117
$ws.Range("Q103").Value = 60

# Row 104
$ws.Range("D104").Value = 44433
$ws.Range("J104").Value = 300
$ws.Range("K104").Value = 7000
$ws.Range("L104").Value = 7000
$ws.Range("M104").Value = 7000
$ws.Range("N104").Value = '$/caja 60 unidades'
$ws.Range("O104").Value = 'Región de Arica y Parinacota'
$ws.Range("P104").Value = 117
$ws.Range("Q104").Value = 60

# Row 105
$ws.Range("D105").Value = 44676
$ws.Range("J105").Value = 150
$ws.Range("K105").Value = 9000
$ws.Range("L105").Value = 9000
$ws.Range("M105").Value = 9000
$ws.Range("N105").Value = '$/caja 50 unidades'
$ws.Range("O105").Value = 'Región del Maule'
$ws.Range("P105").Value = 180
$ws.Range("Q105").Value = 50

# Row 106
$ws.Range("D106").Value = 44417
$ws.Range("J106").Value = 300
$ws.Range("K106").Value = 7000
$ws.Range("L106").Value = 7000
$ws.Range("M106").Value = 7000
$ws.Range("N106").Value = '$/caja 60 unidades'
$ws.Range("O106").Value = 'Región de Arica y Parinacota'
$ws.Range("P106").Value = 117
$ws.Range("Q106").Value = 60

# Row 107
$ws.Range("D107").Value = 44648
$ws.Range("J107").Value = 200
$ws.Range("K107").Value = 7000
$ws.Range("L107").Value = 7000
$ws.Range("M107").Value = 7000
$ws.Range("N107").Value = '$/caja 50 unidades'
$ws.Range("O107").Value = 'Región del Maule'
$ws.Range("P107").Value = 140
$ws.Range("Q107").Value = 50

# Row 108
$ws.Range("D108").Value = 44551
$ws.Range("J108").Value = 100
$ws.Range("K108").Value = 10000
$ws.Range("L108").Value = 10000
$ws.Range("M108").Value = 10000
$ws.Range("N108").Value = '$/caja 50 unidades'
$ws.Range("O108").Value = 'Región del Maule'
$ws.Range("P108").Value = 200
$ws.Range("Q108").Value = 50

# Row 109
$ws.Range("D109").Value = 44554
$ws.Range("J109").Value = 100
$ws.Range("K109").Value = 10000
$ws.Range("L109").Value = 10000
$ws.Range("M109").Value = 10000
$ws.Range("N109").Value = '$/caja 50 unidades'
$ws.Range("O109").Value = 'Región del Maule'
$ws.Range("P109").Value = 200
$ws.Range("Q109").Value = 50

# Row 110
$ws.Range("D110").Value = 44565
$ws.Range("J110").Value = 150
$ws.Range("K110").Value = 9000
$ws.Range("L110").Value = 9000
$ws.Range("M110").Value = 9000
$ws.Range("N110").Value = '$/caja 50 unidades'
$ws.Range("O110").Value = 'Región del Maule'
$ws.Range("P110").Value = 180
$ws.Range("Q110").Value = 50

# Row 111
$ws.Range("D111").Value = 44603
$ws.Range("J111").Value = 200
$ws.Range("K111").Value = 7000
$ws.Range("L111").Value = 7000
$ws.Range("M111").Value = 7000
$ws.Range("N111").Value = '$/caja 50 unidades'
$ws.Range("O111").Value = 'Región del Maule'
$ws.Range("P111").Value = 140
$ws.Range("Q111").Value = 50

# Row 112
$ws.Range("D112").Value = 44263
$ws.Range("J112").Value = 200
$ws.Range("K112").Value = 8000
$ws.Range("L112").Value = 8000
$ws.Range("M112").Value = 8000
$ws.Range("N112").Value = '$/caja 60 unidades'
$ws.Range("O112").Value = 'Región del Maule'
$ws.Range("P112").Value = 133
$ws.Range("Q112").Value = 60

# Row 113
$ws.Range("D113").Value = 44609
$ws.Range("J113").Value = 150
$ws.Range("K113").Value = 7000
$ws.Range("L113").Value = 7000
$ws.Range("M113").Value = 7000
$ws.Range("N113").Value = '$/caja 50 unidades'
$ws.Range("O113").Value = 'Región del Maule'
$ws.Range("P113").Value = 140
$ws.Range("Q113").Value = 50

# Row 114
$ws.Range("D114").Value = 44277
$ws.Range("J114").Value = 200
$ws.Range("K114").Value = 8000
$ws.Range("L114").Value = 8000
$ws.Range("M114").Value = 8000
$ws.Range("N114").Value = '$/caja 60 unidades'
$ws.Range("O114").Value = 'Región del Maule'
$ws.Range("P114").Value = 133
$ws.Range("Q114").Value = 60

# Row 115
$ws.Range("D115").Value = 44265
$ws.Range("J115").Value = 200
$ws.Range("K115").Value = 8000
$ws.Range("L115").Value = 8000
$ws.Range("M115").Value = 8000
$ws.Range("N115").Value = '$/caja 60 unidades'
$ws.Range("O115").Value = 'Región del Maule'
$ws.Range("P115").Value = 133
$ws.Range("Q115").Value = 60

# Row 116
$ws.Range("D116").Value = 44627
$ws.Range("J116").Value = 150
$ws.Range("K116").Value = 7000
$ws.Range("L116").Value = 7000
$ws.Range("M116").Value = 7000
$ws.Range("N116").Value = '$/caja 50 unidades'
$ws.Range("O116").Value = 'Región del Maule'
$ws.Range("P116").Value = 140
$ws.Range("Q116").Value = 50

# Row 117
$ws.Range("D117").Value = 44245
$ws.Range("J117").Value = 200
$ws.Range("K117").Value = 8000
$ws.Range("L117").Value = 8000
$ws.Range("M117").Value = 8000
$ws.Range("N117").Value = '$/caja 60 unidades'
$ws.Range("O117").Value = 'Región del Maule'
$ws.Range("P117").Value = 133
$ws.Range("Q117").Value = 60

# Row 118
$ws.Range("D118").Value = 44249
$ws.Range("J118").Value = 200
$ws.Range("K118").Value = 8000
$ws.Range("L118").Value = 8000
$ws.Range("M118").Value = 8000
$ws.Range("N118").Value = '$/caja 60 unidades'
$ws.Range("O118").Value = 'Región del Maule'
$ws.Range("P118").Value = 133
$ws.Range("Q118").Value = 60

# Row 119
$ws.Range("D119").Value = 44431
$ws.Range("J119").Value = 400
$ws.Range("K119").Value = 6000
$ws.Range("L119").Value = 6000
$ws.Range("M119").Value = 6000
$ws.Range("N119").Value = '$/caja 60 unidades'
$ws.Range("O119").Value = 'Región de Arica y Parinacota'
$ws.Range("P119").Value = 100
$ws.Range("Q119").Value = 60

# Row 120
$ws.Range("D120").Value = 44299
$ws.Range("J120").Value = 200
$ws.Range("K120").Value = 8000
$ws.Range("L120").Value = 8000
$ws.Range("M120").Value = 8000
$ws.Range("N120").Value = '$/caja 60 unidades'
$ws.Range("O120").Value = 'Región del Maule'
$ws.Range("P120").Value = 133
$ws.Range("Q120").Value = 60

# Row 121
$ws.Range("D121").Value = 44615
$ws.Range("J121").Value = 300
$ws.Range("K121").Value = 5000
$ws.Range("L121").Value = 5000
$ws.Range("M121").Value = 5000
$ws.Range("N121").Value = '$/caja 50 unidades'
$ws.Range("O121").Value = 'Región de Arica y Parinacota'
$ws.Range("P121").Value = 100
$ws.Range("Q121").Value = 50

# Row 122
$ws.Range("D122").Value = 44615
$ws.Range("J122").Value = 200
$ws.Range("K122").Value = 6000
$ws.Range("L122").Value = 6000
$ws.Range("M122").Value = 6000
$ws.Range("N122").Value = '$/caja 50 unidades'
$ws.Range("O122").Value = 'Región del Maule'
$ws.Range("P122").Value = 120
$ws.Range("Q122").Value = 50

# Row 123
$ws.Range("D123").Value = 44453
$ws.Range("J123").Value = 200
$ws.Range("K123").Value = 7000
$ws.Range("L123").Value = 7000
$ws.Range("M123").Value = 7000
$ws.Range("N123").Value = '$/caja 50 unidades'
$ws.Range("O123").Value = 'Región de Arica y Parinacota'
$ws.Range("P123").Value = 140
$ws.Range("Q123").Value = 50

# Row 124
$ws.Range("D124").Value = 44421
$ws.Range("J124").Value = 200
$ws.Range("K124").Value = 8000
$ws.Range("L124").Value = 8000
$ws.Range("M124").Value = 8000
$ws.Range("N124").Value = '$/caja 60 unidades'
$ws.Range("O124").Value = 'Región de Arica y Parinacota'
$ws.Range("P124").Value = 133
$ws.Range("Q124").Value = 60

# Row 125
$ws.Range("D125").Value = 44251
$ws.Range("J125").Value = 200
$ws.Range("K125").Value = 8000
$ws.Range("L125").Value = 8000
$ws.Range("M125").Value = 8000
$ws.Range("N125").Value = '$/caja 60 unidades'
$ws.Range("O125").Value = 'Región del Maule'
$ws.Range("P125").Value = 133
$ws.Range("Q125").Value = 60

# Row 126
$ws.Range("D126").Value = 44586
$ws.Range("J126").Value = 200
$ws.Range("K126").Value = 7000
$ws.Range("L126").Value = 7000
$ws.Range("M126").Value = 7000
$ws.Range("N126").Value = '$/caja 50 unidades'
$ws.Range("O126").Value = 'Región del Maule'
$ws.Range("P126").Value = 140
$ws.Range("Q126").Value = 50

# Row 127
$ws.Range("D127").Value = 44617
$ws.Range("J127").Value = 200
$ws.Range("K127").Value = 5000
$ws.Range("L127").Value = 5000
$ws.Range("M127").Value = 5000
$ws.Range("N127").Value = '$/caja 50 unidades'
$ws.Range("O127").Value = 'Región de Arica y Parinacota'
$ws.Range("P127").Value = 100
$ws.Range("Q127").Value = 50

# Row 128
$ws.Range("D128").Value = 44617
$ws.Range("J128").Value = 150
$ws.Range("K128").Value = 6000
$ws.Range("L128").Value = 6000
$ws.Range("M128").Value = 6000
$ws.Range("N128").Value = '$/caja 50 unidades'
$ws.Range("O128").Value = 'Región del Maule'
$ws.Range("P128").Value = 120
$ws.Range("Q128").Value = 50

# Row 129
$ws.Range("D129").Value = 44567
$ws.Range("J129").Value = 300
$ws.Range("K129").Value = 9000
$ws.Range("L129").Value = 9000
$ws.Range("M129").Value = 9000
$ws.Range("N129").Value = '$/caja 50 unidades'
$ws.Range("O129").Value = 'Región del Maule'
$ws.Range("P129").Value = 180
$ws.Range("Q129").Value = 50
$ws.Range("A129").Value = 5
$ws.Range("B129").Value = 'Macroferia Regional de Talca'
$ws.Range("C129").Value = 'Maule'
$ws.Range("E129").Value = 7
$ws.Range("F129").Value = 100112001
$ws.Range("G129").Value = 'Berenjena'
$ws.Range("H129").Value = 'Sin especificar'
$ws.Range("I129").Value = 'Primera'
$ws.Range("R129").Value = 'Hortaliza'
$ws.Range("D129").NumberFormat = "YYYY-MM-DD HH:MM:SS"
